# Update "Forecast Comparison" sheet with a new Week_Start_Date column and
# corrected week labels / is_holiday_week boolean typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1) Insert a new column before the ASIN column (column B), shifting
#    ASIN..is_holiday_week one column to the right (B:I -> C:J).
$ws.Columns.Item(2).Insert()

# 2) Header for the newly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# 3) Week labels (col A) lose their leading zero, and the new
#    Week_Start_Date column (col B) gets the week's start date as plain
#    text (not an Excel date serial), so format the column as Text first.
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

# 4) is_holiday_week (now column J) should hold boolean FALSE values
#    instead of numeric 0.
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}
